# Refreshes the cryptocurrency Price/Volume(1h) columns to the latest scrape
# (GitHub Actions run on Wed Nov 22 10:33:32 UTC 2023), and re-orders the
# last three rows (FTXToken now outranks MXToken; SynthetixNetwork swapped
# out for FraxShare).
#
# Note: several "Price" values look like plain decimals (e.g. "3.84"), which
# Excel would otherwise auto-convert to a Number. Those are written with a
# leading apostrophe to force Text, then the cell style is reset back to
# "Normal" so no stray number formatting is left on the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '36.595.72'
$ws.Range("E2").Value = '  -2.11%  '
$ws.Range("D3").Value = '2.005.88'
$ws.Range("E3").Value = '  -0.41%  '
$ws.Range("E4").Value = '  +0.06%  '
$ws.Range("D5").Value = '''234.72'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -9.84%  '
$ws.Range("E6").Value = '  -3.17%  '
$ws.Range("D8").Value = '''54.84'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -4.05%  '
$ws.Range("D9").Value = '''0.370'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -3.63%  '
$ws.Range("D10").Value = '''57.21'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.52%  '
$ws.Range("D11").Value = '''0.0746'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -4.02%  '
$ws.Range("D12").Value = '''0.0983'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -3.71%  '
$ws.Range("D13").Value = '2.301.30'
$ws.Range("E13").Value = '  -0.34%  '
$ws.Range("D14").Value = '''14.13'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.32%  '
$ws.Range("D15").Value = '''20.16'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -6.66%  '
$ws.Range("D16").Value = '''0.756'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -5.52%  '
$ws.Range("D17").Value = '''5.09'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -2.76%  '
$ws.Range("D18").Value = '2.030.28'
$ws.Range("E18").Value = '  +1.38%  '
$ws.Range("D19").Value = '36.932.05'
$ws.Range("E19").Value = '  -0.93%  '
$ws.Range("D20").Value = '''67.70'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -3.17%  '
$ws.Range("D21").Value = '0.0₃0794'
$ws.Range("E21").Value = '  -5.24%  '
$ws.Range("D22").Value = '''5.32'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +3.86%  '
$ws.Range("D23").Value = '''221.47'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -4.76%  '
$ws.Range("E24").Value = '  +0.00%  '
$ws.Range("E25").Value = '  +2.58%  '
$ws.Range("D26").Value = '''2.39'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -8.12%  '
$ws.Range("D27").Value = '''163.26'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.13%  '
$ws.Range("D28").Value = '''8.64'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -3.60%  '
$ws.Range("D29").Value = '''1.38'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +4.13%  '
$ws.Range("D30").Value = '''0.127'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -1.47%  '
$ws.Range("D31").Value = '''18.73'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -4.43%  '
$ws.Range("E32").Value = '  -2.52%  '
$ws.Range("D33").Value = '''4.35'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -5.60%  '
$ws.Range("D34").Value = '''0.0604'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -6.06%  '
$ws.Range("D35").Value = '''2.40'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +1.36%  '
$ws.Range("D36").Value = '''4.24'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -6.38%  '
$ws.Range("E37").Value = '  -0.12%  '
$ws.Range("D38").Value = '''3.33'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.91%  '
$ws.Range("E39").Value = '  -2.39%  '
$ws.Range("D40").Value = '''5.72'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +4.13%  '
$ws.Range("E41").Value = '  -2.20%  '
$ws.Range("D42").Value = '1.460.66'
$ws.Range("E42").Value = '  +2.40%  '
$ws.Range("D43").Value = '''0.0923'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.63%  '
$ws.Range("D45").Value = '''89.74'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.15%  '
$ws.Range("D46").Value = '''1.10'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -7.95%  '
$ws.Range("D47").Value = '''15.34'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -2.75%  '
$ws.Range("D48").Value = '''0.999'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -3.24%  '
$ws.Range("B49").Value = 'FTXToken'
$ws.Range("C49").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D49").Value = '''3.84'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +25.56%  '
$ws.Range("B50").Value = 'MXToken'
$ws.Range("C50").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D50").Value = '''2.88'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.84%  '
$ws.Range("B51").Value = 'FraxShare'
$ws.Range("C51").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D51").Value = '''6.85'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -2.54%  '
